$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Exact "76÷7=10, 6" "28÷9=3, 1"
Replace-Exact "79÷4=19, 3" "81÷8=10, 1"
Replace-Exact "92÷3=30, 2" "31÷5=6, 1"
Replace-Exact "86÷3=28, 2" "74÷7=10, 4"
Replace-Exact "39÷7=5, 4" "48÷5=9, 3"
Replace-Exact "36÷2=18, 0" "96÷3=32, 0"
Replace-Exact "30÷7=4, 2" "68÷2=34, 0"
Replace-Exact "13÷5=2, 3" "58÷4=14, 2"
Replace-Exact "78÷3=26, 0" "41÷4=10, 1"
Replace-Exact "47÷7=6, 5" "23÷5=4, 3"
Replace-Exact "20÷9=2, 2" "55÷4=13, 3"
Replace-Exact "99÷8=12, 3" "50÷3=16, 2"
Replace-Exact "57÷3=19, 0" "51÷7=7, 2"
Replace-Exact "62÷9=6, 8" "71÷5=14, 1"
Replace-Exact "31÷6=5, 1" "43÷4=10, 3"
Replace-Exact "77÷6=12, 5" "48÷2=24, 0"
Replace-Exact "69÷8=8, 5" "23÷9=2, 5"
Replace-Exact "23÷7=3, 2" "46÷4=11, 2"
Replace-Exact "71÷7=10, 1" "79÷7=11, 2"
Replace-Exact "55÷8=6, 7" "65÷4=16, 1"
Replace-Exact "58÷2=29, 0" "23÷8=2, 7"
Replace-Exact "64÷8=8, 0" "16÷4=4, 0"
Replace-Exact "18÷6=3, 0" "38÷5=7, 3"
Replace-Exact "66÷4=16, 2" "77÷6=12, 5"
Replace-Exact "55÷2=27, 1" "52÷6=8, 4"
